# Apply updated cryptocurrency price/volume data as per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.474.92"
$ws.Range("E2").Value = "  -2.46%  "

# Row 3
$ws.Range("D3").Value = "'1.824.17"
$ws.Range("E3").Value = "  -2.82%  "

# Row 4
$ws.Range("E4").Value = "  +0.57%  "

# Row 5
$ws.Range("D5").Value = "'312.39"
$ws.Range("E5").Value = "  -0.96%  "

# Row 6
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.52%  "

# Row 7
$ws.Range("D7").Value = "'0.4240"
$ws.Range("E7").Value = "  -1.72%  "

# Row 8
$ws.Range("D8").Value = "'0.3585"
$ws.Range("E8").Value = "  -2.82%  "

# Row 9
$ws.Range("D9").Value = "'0.07182"
$ws.Range("E9").Value = "  -3.23%  "

# Row 10
$ws.Range("D10").Value = "'0.8559"
$ws.Range("E10").Value = "  -4.77%  "

# Row 11
$ws.Range("D11").Value = "'20.43"
$ws.Range("E11").Value = "  -3.50%  "

# Row 12
$ws.Range("D12").Value = "'1.761.90"
$ws.Range("E12").Value = "  -3.90%  "

# Row 13
$ws.Range("D13").Value = "'5.345"
$ws.Range("E13").Value = "  -1.69%  "

# Row 14
$ws.Range("D14").Value = "'6.440"
$ws.Range("E14").Value = "  -3.27%  "

# Row 15
$ws.Range("D15").Value = "'0.06925"
$ws.Range("E15").Value = "  -1.34%  "

# Row 16
$ws.Range("D16").Value = "'1.005"
$ws.Range("E16").Value = "  +0.55%  "

# Row 17
$ws.Range("D17").Value = "'81.35"
$ws.Range("E17").Value = "  -0.20%  "

# Row 18
$ws.Range("D18").Value = "'0.000008878"
$ws.Range("E18").Value = "  -2.30%  "

# Row 19
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.38%  "

# Row 20
$ws.Range("D20").Value = "'15.32"
$ws.Range("E20").Value = "  -2.02%  "

# Row 21
$ws.Range("D21").Value = "'27.292.54"
$ws.Range("E21").Value = "  -2.57%  "

# Row 22
$ws.Range("D22").Value = "'5.119"
$ws.Range("E22").Value = "  +0.78%  "

# Row 23
$ws.Range("D23").Value = "'10.87"
$ws.Range("E23").Value = "  +1.24%  "

# Row 24
$ws.Range("D24").Value = "'2.021.28"
$ws.Range("E24").Value = "  -2.23%  "

# Row 25
$ws.Range("D25").Value = "'1.988"
$ws.Range("E25").Value = "  -0.10%  "

# Row 26
$ws.Range("D26").Value = "'154.39"
$ws.Range("E26").Value = "  -0.17%  "

# Row 27
$ws.Range("D27").Value = "'18.48"
$ws.Range("E27").Value = "  -1.69%  "

# Row 28
$ws.Range("D28").Value = "'5.115"
$ws.Range("E28").Value = "  -5.87%  "

# Row 29
$ws.Range("D29").Value = "'113.89"
$ws.Range("E29").Value = "  -4.79%  "

# Row 30
$ws.Range("D30").Value = "'1.772"
$ws.Range("E30").Value = "  -8.40%  "

# Row 31
$ws.Range("D31").Value = "'0.08882"
$ws.Range("E31").Value = "  -1.05%  "

# Row 32
$ws.Range("D32").Value = "'0.7426"
$ws.Range("E32").Value = "  -7.20%  "

# Row 33
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'2.963"
$ws.Range("E33").Value = "  -1.93%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.523"
$ws.Range("E34").Value = "  -3.51%  "

# Row 35
$ws.Range("D35").Value = "'1.119"
$ws.Range("E35").Value = "  -6.32%  "

# Row 36
$ws.Range("D36").Value = "'1.002"
$ws.Range("E36").Value = "  +0.52%  "

# Row 37
$ws.Range("D37").Value = "'1.079"
$ws.Range("E37").Value = "  -4.26%  "

# Row 38
$ws.Range("D38").Value = "'0.05232"
$ws.Range("E38").Value = "  -4.92%  "

# Row 39
$ws.Range("D39").Value = "'0.01919"
$ws.Range("E39").Value = "  -2.17%  "

# Row 40
$ws.Range("D40").Value = "'2.784"
$ws.Range("E40").Value = "  -3.25%  "

# Row 41
$ws.Range("D41").Value = "'0.5056"
$ws.Range("E41").Value = "  -2.44%  "

# Row 42
$ws.Range("D42").Value = "'0.1651"
$ws.Range("E42").Value = "  -2.91%  "

# Row 43
$ws.Range("D43").Value = "'6.339"
$ws.Range("E43").Value = "  -7.17%  "

# Row 44
$ws.Range("D44").Value = "'8.319"
$ws.Range("E44").Value = "  -3.52%  "

# Row 45
$ws.Range("D45").Value = "'10.46"
$ws.Range("E45").Value = "  -1.58%  "

# Row 46
$ws.Range("D46").Value = "'106.10"
$ws.Range("E46").Value = "  -0.44%  "

# Row 47
$ws.Range("D47").Value = "'0.06453"
$ws.Range("E47").Value = "  -2.57%  "

# Row 48
$ws.Range("D48").Value = "'0.4665"
$ws.Range("E48").Value = "  -2.48%  "

# Row 49
$ws.Range("E49").Value = "  +0.51%  "

# Row 50
$ws.Range("D50").Value = "'1.608"
$ws.Range("E50").Value = "  -3.14%  "

# Row 51
$ws.Range("D51").Value = "'63.50"
$ws.Range("E51").Value = "  -3.03%  "
